# "switch to tau from mu" - rename observed-timepoint column headers
# (formerly growth-rate mu labels "u_x.y") to elapsed-time tau labels
# "t_N_min", and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "t_100_min"
$ws.Range("E1").Value = "t_60_min"
$ws.Range("F1").Value = "t_40_min"
$ws.Range("G1").Value = "t_30_min"
$ws.Range("H1").Value = "t_24_min"
$ws.Range("I1").Value = "t_20_min"

$ws.Range("H11").Select()
